# "Test mit allen Versuchen" -- record the additional LDAP/Webshell port
# entries discovered while testing the HackerLab services, and correct the
# Log4Shell LDAP listener port.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group "Password_Cracking & SQLi" (rows 9-11): give it the same
#     centered / word-wrapped look as the other merged service groups
#     above it, then merge the label column. ---
$ws.Range("A2").Copy()
$ws.Range("A9:A11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A9:A11").Merge()

# --- New "Webshell" service block (rows 12-13) ---
$ws.Range("A12").Value = "Webshell"
$ws.Range("B12").Value = "database"
$ws.Range("C12").Value = "3307:3307"

$ws.Range("B13").Value = "web-server"
$ws.Range("C13").Value = "8082:82"

# Match the (non-wrapped) centered look used by the "Log4Shell_Vulnerable"
# group and merge the new label column.
$ws.Range("A5").Copy()
$ws.Range("A12:A13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12:A13").Merge()

$excel.CutCopyMode = $false

# --- Correct the LDAP port used by the patched Log4Shell service ---
$ws.Range("C8").Value = "10389:10390"

# --- Restore the cursor position saved with the workbook ---
$ws.Range("F7").Select()
